$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 120
$ws.Range("AI2").Value = 55
$ws.Range("AO2").Value = 46
$ws.Range("W2").Value = 2.1
$ws.Range("Z2").Value = 48
$ws.Range("AF3").Value = 9
$ws.Range("AN3").Value = 5.6
$ws.Range("H3").Value = 8.4
$ws.Range("J3").Value = 5.1
$ws.Range("Q3").Value = 1.63
$ws.Range("S3").Value = 2.58
$ws.Range("U3").Value = 2.08
$ws.Range("V3").Value = 1.12
$ws.Range("W3").Value = 3.15
$ws.Range("F4").Value = 2.02
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 3.75
$ws.Range("J4").Value = 4.1
$ws.Range("N4").Value = 5.1
$ws.Range("S4").Value = 2.64
$ws.Range("V4").Value = 1.36
$ws.Range("W4").Value = 1.96
$ws.Range("F5").Value = 1.25
$ws.Range("G5").Value = 1.87
$ws.Range("H5").Value = 4.1
$ws.Range("I5").Value = 7.4
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 950
$ws.Range("L5").Value = 1.17
$ws.Range("N5").Value = 2.62
$ws.Range("O5").Value = 1.08
$ws.Range("P5").Value = 2.62
$ws.Range("Q5").Value = 1.08
$ws.Range("R5").Value = 1.96
$ws.Range("S5").Value = 1.61
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.04
$ws.Range("V5").Value = 1.18
$ws.Range("W5").Value = 2.14
$ws.Range("F6").Value = 2.34
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 2.7
$ws.Range("N6").Value = 2.18
$ws.Range("P6").Value = 2.16
$ws.Range("R6").Value = 1.48
$ws.Range("S6").Value = 2.24
$ws.Range("T6").Value = 1.04
$ws.Range("V6").Value = 1.5
$ws.Range("F7").Value = 1.43
$ws.Range("G7").Value = 1.72
$ws.Range("H7").Value = 5.3
$ws.Range("I7").Value = 10.5
$ws.Range("J7").Value = 3.85
$ws.Range("L7").Value = 1.22
$ws.Range("N7").Value = 2.26
$ws.Range("R7").Value = 1.35
$ws.Range("S7").Value = 2.28
$ws.Range("T7").Value = 1.04
$ws.Range("U7").Value = 1.04
$ws.Range("V7").Value = 1.1
$ws.Range("W7").Value = 2.38
$ws.Range("X7").Value = 26
$ws.Range("AL8").Value = 44
$ws.Range("AM8").Value = 110
$ws.Range("G8").Value = 2.48
$ws.Range("H8").Value = 3.35
$ws.Range("P8").Value = 1.81
$ws.Range("AJ9").Value = 38
$ws.Range("AO9").Value = 12.5
$ws.Range("F9").Value = 2.64
$ws.Range("G9").Value = 2.68
$ws.Range("I9").Value = 2.66
$ws.Range("V9").Value = 1.6
$ws.Range("W9").Value = 1.59
$ws.Range("AF10").Value = 90
$ws.Range("AK10").Value = 140
$ws.Range("AM10").Value = 130
$ws.Range("AO10").Value = 5.1
$ws.Range("H10").Value = 1.38
$ws.Range("I10").Value = 1.39
$ws.Range("J10").Value = 5.6
$ws.Range("K10").Value = 5.8
$ws.Range("P10").Value = 2.5
$ws.Range("T10").Value = 1.95
$ws.Range("V10").Value = 3.55
$ws.Range("AC11").Value = 21
$ws.Range("AE11").Value = 470
$ws.Range("AF11").Value = 8
$ws.Range("AG11").Value = 13.5
$ws.Range("AH11").Value = 130
$ws.Range("AI11").Value = 310
$ws.Range("H11").Value = 20
$ws.Range("I11").Value = 22
$ws.Range("P11").Value = 2.88
$ws.Range("Q11").Value = 1.48
$ws.Range("R11").Value = 1.78
$ws.Range("X11").Value = 36
$ws.Range("AB12").Value = 15.5
$ws.Range("AI12").Value = 110
$ws.Range("AO12").Value = 100
$ws.Range("J12").Value = 6.8
$ws.Range("P12").Value = 3.4
$ws.Range("U12").Value = 2.28
$ws.Range("X12").Value = 40
$ws.Range("AA13").Value = 15.5
$ws.Range("AN13").Value = 90
$ws.Range("G13").Value = 6.4
$ws.Range("I13").Value = 1.64
$ws.Range("J13").Value = 4.4
$ws.Range("K13").Value = 4.5
$ws.Range("Q13").Value = 1.83
$ws.Range("V13").Value = 2.56
$ws.Range("X13").Value = 16
$ws.Range("AA14").Value = 32
$ws.Range("AB14").Value = 15.5
$ws.Range("AH14").Value = 15.5
$ws.Range("AJ14").Value = 1000
$ws.Range("AM14").Value = 65
$ws.Range("AO14").Value = 15
$ws.Range("F14").Value = 3.1
$ws.Range("M14").Value = 1.05
$ws.Range("O14").Value = 1.24
$ws.Range("P14").Value = 2.28
$ws.Range("Q14").Value = 1.75
$ws.Range("R14").Value = 1.51
$ws.Range("S14").Value = 2.86
$ws.Range("X14").Value = 19.5
$ws.Range("Y14").Value = 13.5
$ws.Range("F15").Value = 1.7
$ws.Range("G15").Value = 3.2
$ws.Range("H15").Value = 1.81
$ws.Range("I15").Value = 3.5
$ws.Range("J15").Value = 2.48
$ws.Range("K15").Value = 950
$ws.Range("L15").Value = 1.01
$ws.Range("S15").Value = 2.24
$ws.Range("W15").Value = 1.45
$ws.Range("AA16").Value = 65
$ws.Range("AB16").Value = 12.5
$ws.Range("AC16").Value = 9.199999999999999
$ws.Range("AD16").Value = 16.5
$ws.Range("AE16").Value = 46
$ws.Range("AF16").Value = 21
$ws.Range("AG16").Value = 15
$ws.Range("AH16").Value = 22
$ws.Range("AI16").Value = 60
$ws.Range("AJ16").Value = 48
$ws.Range("AK16").Value = 38
$ws.Range("AL16").Value = 55
$ws.Range("AM16").Value = 130
$ws.Range("AN16").Value = 34
$ws.Range("AO16").Value = 44
$ws.Range("F16").Value = 2.48
$ws.Range("G16").Value = 2.76
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 3.25
$ws.Range("R16").Value = 1.3
$ws.Range("S16").Value = 3.6
$ws.Range("T16").Value = 1.77
$ws.Range("U16").Value = 2.04
$ws.Range("W16").Value = 1.58
$ws.Range("X16").Value = 15.5
$ws.Range("Y16").Value = 14
$ws.Range("Z16").Value = 24
